$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the full data (columns B:AC) between paired rows. The "id" column
#    (A) keeps its sequential value in each row; everything else (match
#    id, date, teams, odds, PL columns) is exchanged between the two rows.
# ---------------------------------------------------------------------------
function Swap-Rows($ws, $r1, $r2) {
    $rng1 = $ws.Range("B$r1" + ":AC$r1")
    $rng2 = $ws.Range("B$r2" + ":AC$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

Swap-Rows $ws 30 31
Swap-Rows $ws 177 178
Swap-Rows $ws 204 205

# ---------------------------------------------------------------------------
# 2) Refresh the trailing fixtures block. The three still-upcoming fixtures
#    (previously rows 269-271) move up to rows 264-266 with refreshed odds,
#    and the now-stale rows 267-271 are removed entirely.
# ---------------------------------------------------------------------------

function Set-RowValues($ws, $row, $values) {
    # $values is an ordered hashtable-like array of (column, value) pairs
    foreach ($pair in $values) {
        $ws.Range($pair[0] + $row).Value2 = $pair[1]
    }
}

# Row 264 <- updated "Casa Pia vs FC Porto" fixture
Set-RowValues $ws 264 @(
    ,@("B", 6961656)
    ,@("E", 45403.58333333334)
    ,@("F", "Casa Pia")
    ,@("G", "FC Porto")
    ,@("K", 8.5)
    ,@("L", 5)
    ,@("M", 1.363)
    ,@("N", 8)
    ,@("O", 5.25)
    ,@("P", 1.333)
    ,@("Q", 1.5)
    ,@("R", 1.88)
    ,@("S", 2.02)
    ,@("T", 2.75)
    ,@("U", 1.9)
    ,@("V", 1.95)
)

# Row 265 <- updated "Sporting vs Guimaraes" fixture
Set-RowValues $ws 265 @(
    ,@("B", 7024015)
    ,@("E", 45403.6875)
    ,@("F", "Sporting")
    ,@("G", "Guimaraes")
    ,@("K", 1.333)
    ,@("L", 5)
    ,@("M", 9)
    ,@("N", 1.285)
    ,@("O", 5.5)
    ,@("P", 11)
    ,@("Q", -1.75)
    ,@("R", 2.04)
    ,@("S", 1.86)
    ,@("T", 3)
    ,@("U", 1.85)
    ,@("V", 2)
)

# Row 266 <- updated "SC Farense vs Benfica" fixture
Set-RowValues $ws 266 @(
    ,@("B", 7023334)
    ,@("E", 45404.67708333334)
    ,@("F", "SC Farense")
    ,@("G", "Benfica")
    ,@("K", 8)
    ,@("L", 4.5)
    ,@("M", 1.4)
    ,@("N", 7)
    ,@("O", 5.25)
    ,@("P", 1.4)
    ,@("Q", 1.25)
    ,@("R", 2.05)
    ,@("S", 1.85)
    ,@("T", 3)
    ,@("U", 1.875)
    ,@("V", 1.975)
)

# Remove the now-obsolete trailing rows (267-271)
$ws.Range("A267:AC271").EntireRow.Delete()
